# Updated tourism data. Includes sum within islands for BES, KNA and TTO
#
# The source dataset had per-sub-island breakout rows for three
# multi-island entries (Bonaire/Sint Eustatius/Saba, Saint Kitts and
# Nevis, Trinidad and Tobago). Those breakout rows are removed and the
# parent row now carries the summed foreign_tourists figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bonaire, Sint Eustatius and Saba --------------------------------
# Row 7 is the "Bonaire, Sint Eustatius and Saba" total row (previously
# "NA"); rows 8-10 are the Bonaire / Sint Eustatius / Saba breakouts
# (133000 + 14600 + 20500 = 168100). Put the summed total on row 7 and
# drop the three breakout rows.
$ws.Range("E7").Value = 168100
$ws.Rows("8:10").Delete()

# --- Saint Kitts and Nevis -------------------------------------------
# After the deletion above, the Saint Kitts / Nevis breakout rows
# (previously "NA", original rows 25-26) are now at rows 22-23. The
# "Saint Kitts and Nevis" total row already carries the correct summed
# value (104730), so just remove the breakout rows.
$ws.Rows("22:23").Delete()

# --- Trinidad and Tobago ----------------------------------------------
# After the two deletions above, the Trinidad / Tobago breakout rows
# (previously "NA", original rows 32-33) are now at rows 27-28. The
# "Trinidad and Tobago" total row already carries the correct summed
# value (412537), so just remove the breakout rows.
$ws.Rows("27:28").Delete()

# Restore the active selection recorded in the saved workbook.
$ws.Range("G18").Select()
